# Updated cryptos list on Sun Sep 10 17:45:06 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures for every coin row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'25.953.88"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.19%  "
$c = $ws.Range("D3")
$c.Value = "'1.622.58"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  +0.52%  "
$c = $ws.Range("D5")
$c.Value = "'213.68"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  +0.53%  "
$ws.Range("E8").Value = "  -2.39%  "
$ws.Range("E9").Value = "  -3.07%  "
$c = $ws.Range("D10")
$c.Value = "'18.13"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -7.47%  "
$c = $ws.Range("D11")
$c.Value = "'0.0787"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.89%  "
$c = $ws.Range("D12")
$c.Value = "'1.848.50"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.13%  "
$c = $ws.Range("D13")
$c.Value = "'1.639.25"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("E15").Value = "  -3.76%  "
$c = $ws.Range("D16")
$c.Value = "'25.937.47"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.60%  "
$c = $ws.Range("D17")
$c.Value = "'0.0₃0736"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -3.38%  "
$c = $ws.Range("D18")
$c.Value = "'61.15"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -3.51%  "
$ws.Range("E19").Value = "  +0.51%  "
$c = $ws.Range("D20")
$c.Value = "'189.99"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("E21").Value = "  -2.84%  "
$c = $ws.Range("D22")
$c.Value = "'9.55"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.69%  "
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("E24").Value = "  +0.79%  "
$c = $ws.Range("D25")
$c.Value = "'143.34"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("E28").Value = "  -2.50%  "
$c = $ws.Range("D29")
$c.Value = "'15.15"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.38%  "
$ws.Range("E30").Value = "  -1.68%  "
$ws.Range("E31").Value = "  -3.14%  "
$ws.Range("E32").Value = "  -4.43%  "
$ws.Range("E33").Value = "  -5.76%  "
$c = $ws.Range("D35")
$c.Value = "'1.48"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -3.03%  "
$c = $ws.Range("D36")
$c.Value = "'1.125.85"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.19%  "
$c = $ws.Range("D37")
$c.Value = "'0.842"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -6.87%  "
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("E39").Value = "  -4.88%  "
$ws.Range("E40").Value = "  -2.15%  "
$c = $ws.Range("D41")
$c.Value = "'97.53"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.79%  "
$c = $ws.Range("D42")
$c.Value = "'0.770"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.12%  "
$c = $ws.Range("D43")
$c.Value = "'1.759.31"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("E44").Value = "  -5.45%  "
$ws.Range("E45").Value = "  -2.98%  "
$c = $ws.Range("D46")
$c.Value = "'54.41"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("E48").Value = "  +0.08%  "
$c = $ws.Range("D49")
$c.Value = "'0.413"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("E50").Value = "  +0.57%  "
$c = $ws.Range("D51")
$c.Value = "'7.46"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -3.38%  "
